$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 124, pushing every
# subsequent record (old rows 124-239) down by one row (new rows 125-240).
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row with its data.
$ws.Cells.Item(124, 1).Value  = 10
$ws.Cells.Item(124, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(124, 3).Value  = "La Araucanía"
$ws.Cells.Item(124, 4).Value  = 44447
$ws.Cells.Item(124, 5).Value  = 9
$ws.Cells.Item(124, 6).Value  = 100112023
$ws.Cells.Item(124, 7).Value  = "Brócoli"
$ws.Cells.Item(124, 8).Value  = "Sin especificar"
$ws.Cells.Item(124, 9).Value  = "Primera"
$ws.Cells.Item(124, 10).Value = 1900
$ws.Cells.Item(124, 11).Value = 700
$ws.Cells.Item(124, 12).Value = 800
$ws.Cells.Item(124, 13).Value = 766
$ws.Cells.Item(124, 14).Value = "`$/unidad"
$ws.Cells.Item(124, 15).Value = "Región Metropolitana"
$ws.Cells.Item(124, 16).Value = 766
$ws.Cells.Item(124, 17).Value = 1
$ws.Cells.Item(124, 18).Value = "Hortaliza"
